$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address value in cell B2 (preserves hyperlink formatting)
$ws.Range("B2").Value = "chappel.mann+stl2@gmail.com"

# Update the active cell selection to B2
$ws.Range("B2").Select()
